$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92 (shifts existing rows 92:130 down to 93:131)
$ws.Rows("92:92").Insert()

# Populate the newly inserted row 92 with the new record
$ws.Cells.Item(92, 1).Value = 7
$ws.Cells.Item(92, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(92, 3).Value = "Ñuble"
$ws.Cells.Item(92, 4).Value = 44609
$ws.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(92, 5).Value = 16
$ws.Cells.Item(92, 6).Value = 100112045
$ws.Cells.Item(92, 7).Value = "Zapallo"
$ws.Cells.Item(92, 8).Value = "Camote"
$ws.Cells.Item(92, 9).Value = "1a (cosecha)"
$ws.Cells.Item(92, 10).Value = 200
$ws.Cells.Item(92, 11).Value = 350
$ws.Cells.Item(92, 12).Value = 400
$ws.Cells.Item(92, 13).Value = 375
$ws.Cells.Item(92, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(92, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(92, 16).Value = 375
$ws.Cells.Item(92, 17).Value = 1
$ws.Cells.Item(92, 18).Value = "Hortaliza"
